# Update "想去人数" (F column) values on both the "展览" and "全部类型"
# sheets, which hold identical data tables.

$wb = $excel.ActiveWorkbook

# Map of cell address -> new value to apply on each target worksheet.
$updates = @{
    "F2"  = 2069
    "F12" = 29
    "F14" = 233
    "F15" = 10
    "F19" = 3969
    "F24" = 942
    "F25" = 702
    "F29" = 1770
    "F30" = 32
    "F31" = 32
    "F32" = 63
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($addr in $updates.Keys) {
        $ws.Range($addr).Value = $updates[$addr]
    }
}
